$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = '28.287.99'
$ws.Range("E2").Value = '  +3.18%  '
$ws.Range("D3").Value = '1.816.32'
$ws.Range("E3").Value = '  +4.13%  '
$ws.Range("E4").Value = '  -0.06%  '
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = '328.43'
$ws.Range("E5").Value = '  +2.12%  '
$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = '1.001'
$ws.Range("E6").Value = '  +0.08%  '
$ws.Range("D7").NumberFormat = "@"
$ws.Range("D7").Value = '0.4368'
$ws.Range("E7").Value = '  +3.58%  '
$ws.Range("D8").NumberFormat = "@"
$ws.Range("D8").Value = '0.3672'
$ws.Range("E8").Value = '  +2.46%  '
$ws.Range("D9").NumberFormat = "@"
$ws.Range("D9").Value = '44.97'
$ws.Range("E9").Value = '  -1.18%  '
$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = '0.07701'
$ws.Range("E10").Value = '  +3.84%  '
$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = '1.143'
$ws.Range("E11").Value = '  +2.72%  '
$ws.Range("D12").NumberFormat = "@"
$ws.Range("D12").Value = '1.001'
$ws.Range("E12").Value = '  -0.04%  '
$ws.Range("D13").NumberFormat = "@"
$ws.Range("D13").Value = '22.17'
$ws.Range("E13").Value = '  +3.30%  '
$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = '6.319'
$ws.Range("E14").Value = '  +3.47%  '
$ws.Range("D15").NumberFormat = "@"
$ws.Range("D15").Value = '7.547'
$ws.Range("E15").Value = '  +4.92%  '
$ws.Range("D16").Value = '1.835.25'
$ws.Range("E16").Value = '  +5.52%  '
$ws.Range("D17").NumberFormat = "@"
$ws.Range("D17").Value = '93.00'
$ws.Range("E17").Value = '  +6.18%  '
$ws.Range("E18").Value = '  +1.50%  '
$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = '0.06528'
$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = '1.001'
$ws.Range("E20").Value = '  +0.07%  '
$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = '17.54'
$ws.Range("E21").Value = '  +4.05%  '
$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = '6.272'
$ws.Range("E22").Value = '  +2.77%  '
$ws.Range("D23").Value = '28.322.96'
$ws.Range("E23").Value = '  +3.18%  '
$ws.Range("E24").Value = '  +1.79%  '
$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = '2.040'
$ws.Range("E25").Value = '  -12.85%  '
$ws.Range("D26").NumberFormat = "@"
$ws.Range("D26").Value = '162.43'
$ws.Range("E26").Value = '  +6.78%  '
$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = '20.75'
$ws.Range("E27").Value = '  +1.72%  '
$ws.Range("D28").Value = '2.036.79'
$ws.Range("E28").Value = '  +5.14%  '
$ws.Range("D29").NumberFormat = "@"
$ws.Range("D29").Value = '2.302'
$ws.Range("E29").Value = '  -3.11%  '
$ws.Range("D30").NumberFormat = "@"
$ws.Range("D30").Value = '128.85'
$ws.Range("E30").Value = '  +2.46%  '
$ws.Range("D31").NumberFormat = "@"
$ws.Range("D31").Value = '1.224'
$ws.Range("E31").Value = '  +2.19%  '
$ws.Range("D32").NumberFormat = "@"
$ws.Range("D32").Value = '5.980'
$ws.Range("E32").Value = '  +5.26%  '
$ws.Range("D33").NumberFormat = "@"
$ws.Range("D33").Value = '0.09225'
$ws.Range("E33").Value = '  +1.10%  '
$ws.Range("B34").Value = 'HuobiToken'
$ws.Range("C34").Value = 'https://coinranking.com/coin/DXwP4wF9ksbBO+huobitoken-ht'
$ws.Range("D34").NumberFormat = "@"
$ws.Range("D34").Value = '3.461'
$ws.Range("E34").Value = '  -4.54%  '
$ws.Range("B35").Value = 'Aptos'
$ws.Range("C35").Value = 'https://coinranking.com/coin/HGYj5JCv5+aptos-apt'
$ws.Range("D35").NumberFormat = "@"
$ws.Range("D35").Value = '13.00'
$ws.Range("E35").Value = '  +2.52%  '
$ws.Range("D36").NumberFormat = "@"
$ws.Range("D36").Value = '0.02360'
$ws.Range("E36").Value = '  +2.93%  '
$ws.Range("D37").NumberFormat = "@"
$ws.Range("D37").Value = '0.2185'
$ws.Range("E37").Value = '  +2.00%  '
$ws.Range("D38").NumberFormat = "@"
$ws.Range("D38").Value = '5.200'
$ws.Range("E38").Value = '  +2.23%  '
$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = '0.6597'
$ws.Range("E39").Value = '  +3.26%  '
$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = '0.06212'
$ws.Range("E40").Value = '  +2.64%  '
$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = '8.158'
$ws.Range("E41").Value = '  +3.12%  '
$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = '1.196'
$ws.Range("E42").Value = '  +0.55%  '
$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = '1.437'
$ws.Range("E43").Value = '  +1.26%  '
$ws.Range("E44").Value = '  +0.06%  '
$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = '13.97'
$ws.Range("E45").Value = '  +1.67%  '
$ws.Range("E46").Value = '  +4.61%  '
$ws.Range("E47").Value = '  +1.17%  '
$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = '126.06'
$ws.Range("E48").Value = '  +0.65%  '
$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = '2.027'
$ws.Range("E49").Value = '  +4.08%  '
$ws.Range("E50").Value = '  +4.71%  '
$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = '0.07012'
$ws.Range("E51").Value = '  +2.75%  '
